# Add data for 2022-05-21 (new carjacking counts)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-21"

# Update the "May (through ...)" label in column A, row 6
$ws.Range("A6").Value = "May (through 05-21)"

# Update the May row (row 6) values for years 2016-2022 (columns C-I)
$ws.Range("C6").Value = 33
$ws.Range("E6").Value = 31
$ws.Range("F6").Value = 29
$ws.Range("G6").Value = 42
$ws.Range("H6").Value = 80
$ws.Range("I6").Value = 76

# Update the Total row (row 7) values for years 2016-2022 (columns C-I)
$ws.Range("C7").Value = 195
$ws.Range("E7").Value = 277
$ws.Range("F7").Value = 184
$ws.Range("G7").Value = 304
$ws.Range("H7").Value = 603
$ws.Range("I7").Value = 628
